$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pokemon")
$ws.Activate()

# Insert a new column before column H ("exp apporte" shifts from H to I)
$ws.Columns("H").Insert()

# New "vitesse" column header and values
$ws.Range("H1").Value = "vitesse"
$ws.Range("H2").Value = 10
$ws.Range("H3").Value = 1000
$ws.Range("H4").Value = 10

$ws.Range("H3").Select()
